$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.852.47"
$ws.Range("E2").Value = "  -2.22%  "
$ws.Range("D3").Value = "3.123.03"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.10%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "589.00"
$ws.Range("E5").Value = "  -2.48%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "135.51"
$ws.Range("E6").Value = "  -5.31%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "3.115.92"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("E10").Value = "  -4.44%  "
$ws.Range("E11").Value = "  -2.75%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.452"
$ws.Range("E12").Value = "  -3.42%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000241"
$ws.Range("E13").Value = "  -5.89%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "33.83"
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("D15").Value = "3.633.04"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "62.940.53"
$ws.Range("E17").Value = "  -2.18%  "
$ws.Range("D18").Value = "3.120.94"
$ws.Range("E18").Value = "  -0.22%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.61"
$ws.Range("E19").Value = "  -4.32%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "468.13"
$ws.Range("E20").Value = "  -2.51%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.02"
$ws.Range("E21").Value = "  -3.83%  "
$ws.Range("E22").Value = "  -2.62%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.62"
$ws.Range("E23").Value = "  -1.06%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "85.05"
$ws.Range("E24").Value = "  -0.58%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "12.84"
$ws.Range("E25").Value = "  -4.35%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -2.22%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "7.81"
$ws.Range("E28").Value = "  -6.71%  "
$ws.Range("E29").Value = "  +1.30%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "6.79"
$ws.Range("E30").Value = "  -5.23%  "
$ws.Range("E31").Value = "  +0.12%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "26.49"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("E33").Value = "  -5.43%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "2.52"
$ws.Range("E34").Value = "  -4.92%  "
$ws.Range("E35").Value = "  -3.67%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "51.90"
$ws.Range("E36").Value = "  -0.95%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "5.70"
$ws.Range("E37").Value = "  -4.46%  "
$ws.Range("D38").Value = "0.0₃0673"
$ws.Range("E38").Value = "  -12.94%  "
$ws.Range("E39").Value = "  -2.40%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "413.11"
$ws.Range("E40").Value = "  -7.36%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").Value = "2.892.46"
$ws.Range("E43").Value = "  -12.55%  "
$ws.Range("E44").Value = "  -6.60%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.255"
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("E46").Value = "  +0.13%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.08"
$ws.Range("E47").Value = "  -6.35%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "25.20"
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("E50").Value = "  -8.51%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "120.32"
$ws.Range("E51").Value = "  -0.11%  "
